$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 51, shifting rows 51:55 down to 52:56
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new data record
$ws.Range("A51").Value = 11
$ws.Range("B51").Value = "Vega Monumental Concepción"
$ws.Range("C51").Value = "Bíobío"
$ws.Range("D51").Value = 45106
$ws.Range("E51").Value = 8
$ws.Range("F51").Value = 100114007
$ws.Range("G51").Value = "Jengibre"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 50
$ws.Range("K51").Value = 15000
$ws.Range("L51").Value = 16000
$ws.Range("M51").Value = 15600
$ws.Range("N51").Value = "$/caja 13 kilos"
$ws.Range("O51").Value = "Perú"
$ws.Range("P51").Value = 1200
$ws.Range("Q51").Value = 13
$ws.Range("R51").Value = "Hortaliza"
